$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Marks")

$ws.Range("E4").Value = 6
$ws.Range("E5").Value = "Good but there are 2 problems: 1) pufosos are not correctly calculated (results are not the same as in the document); 2) instead of a cuadratic complextiy you can get a O(nlogn) complexity if you sort the elements beforehand (or if you use a PriorityQueue)"

$ws.Range("E5:E12").Select()
